{"js": "// Change HTTP status code mentioned in the report from 422 to 400, and\n// trim the accompanying explanation so it no longer claims the data was\n// \"syntactically correct\" (the vote endpoint now just says the client\n// sent invalid data).\n//\n// Before: \"...service reply with status code 422, as the client provided\n//          syntactically correct but invalid data to be processed by the\n//          service.\"\n// After:  \"...service reply with status code 400, as the client invalid\n//          data to be processed by the service.\"\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) \"status code 422\" -> \"status code 400\"\nconst codeHits = body.search(\"status code 422\", { matchCase: true, matchWholeWord: false });\ncodeHits.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < codeHits.items.length; i++) {\n  codeHits.items[i].insertText(\"status code 400\", \"Replace\");\n}\nawait context.sync();\n\n// The paragraph also carries an (invisible) \"_GoBack\" bookmark sitting right\n// in the middle of the sentence we are about to edit (\"...be processed by\n// th|e service.\", bookmark between \"th\" and \"e\"). A plain text search/replace\n// across that span would silently drop it, so remove it first and re-create\n// it afterwards at the equivalent spot (right after \"client \").\nlet hadGoBack = true;\ntry {\n  const existing = doc.getBookmarkRange(\"_GoBack\");\n  existing.load(\"text\");\n  await context.sync();\n} catch (e) {\n  hadGoBack = false;\n}\nif (hadGoBack) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 2) \" client provided syntactically correct but invalid data to be\n//    processed by the service.\" -> \" client invalid data to be processed\n//    by the service.\"\nconst explanationHits = body.search(\n  \"client provided syntactically correct but invalid data to be processed by the service.\",\n  { matchCase: true, matchWholeWord: false }\n);\nexplanationHits.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < explanationHits.items.length; i++) {\n  explanationHits.items[i].insertText(\n    \"client invalid data to be processed by the service.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 3) Restore the \"_GoBack\" bookmark right after \"client \" (its position\n// relative to the surrounding, still-present text is unchanged).\nif (hadGoBack) {\n  const anchorHits = body.search(\"client \", { matchCase: true });\n  anchorHits.load(\"text\");\n  await context.sync();\n  if (anchorHits.items.length > 0) {\n    const insertionPoint = anchorHits.items[anchorHits.items.length - 1].getRange(\"End\");\n    insertionPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Change HTTP status code mentioned in the report from 422 to 400, and\n# trim the accompanying explanation so it no longer claims the data was\n# \"syntactically correct\" (the vote endpoint now just says the client\n# sent invalid data).\n\n$d = $word.ActiveDocument\n\n# 1) \"status code 422\" -> \"status code 400\"\n$find1 = $d.Content.Find\n$find1.Text = \"status code 422\"\n$find1.Replacement.Text = \"status code 400\"\n$find1.Execute(\n  $find1.Text,   # FindText\n  $false,        # MatchCase\n  $false,        # MatchWholeWord\n  $false,        # MatchWildcards\n  $false,        # MatchSoundsLike\n  $false,        # MatchAllWordForms\n  $true,         # Forward\n  1,             # Wrap (wdFindContinue)\n  $false,        # Format\n  $find1.Replacement.Text,  # ReplaceWith\n  2              # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) \" client provided syntactically correct but invalid data to be\n#    processed by the service.\" -> \" client invalid data to be processed\n#    by the service.\"\n$find2 = $d.Content.Find\n$find2.Text = \"client provided syntactically correct but invalid data to be processed by the service.\"\n$find2.Replacement.Text = \"client invalid data to be processed by the service.\"\n$find2.Execute(\n  $find2.Text,\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $find2.Replacement.Text,\n  2\n) | Out-Null\n"}
